$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Albahaca" series. It belongs
# right before the current row 332, so push that row (and everything below
# it) down by one and populate the freed-up row with the new record.
$ws.Rows.Item(332).Insert()

$ws.Range("A332").Value = 9
$ws.Range("B332").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C332").Value = "Metropolitana"
$ws.Range("D332").Value = 44841
$ws.Range("D332").NumberFormat = $ws.Range("D333").NumberFormat
$ws.Range("E332").Value = 13
$ws.Range("F332").Value = 100112052
$ws.Range("G332").Value = "Albahaca"
$ws.Range("H332").Value = "Sin especificar"
$ws.Range("I332").Value = "Primera"
$ws.Range("J332").Value = 60
$ws.Range("K332").Value = 7000
$ws.Range("L332").Value = 7000
$ws.Range("M332").Value = 7000
$ws.Range("N332").Value = "$/docena de matas"
$ws.Range("O332").Value = "Provincia de Chacabuco"
$ws.Range("P332").Value = 1167
$ws.Range("Q332").Value = 6
$ws.Range("R332").Value = "Hortaliza"
